$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the analog-buffer sensor rows in column J ("x" = shared string already
# used elsewhere in the column) for the new, better layout of analog buffers.
$ws.Range("J12").Value = "x"
$ws.Range("J14").Value = "x"
$ws.Range("J15").Value = "x"
$ws.Range("J17").Value = "x"
$ws.Range("J18").Value = "x"

# Slightly narrow/widen the first two columns to make room for the new layout.
# (The host's ColumnWidth setter snaps to a pixel grid, same as real Excel;
# these inputs land on the grid points closest to the authored widths of
# 29.375 / 9.625 characters.)
$ws.Columns(1).ColumnWidth = 28.5
$ws.Columns(2).ColumnWidth = 8.833333333333334

# The workbook was left with the selection resting on the last cell touched.
$ws.Range("J18").Select()
